$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new "StatQuery" column needs to be inserted between the existing
# column A (query) and column B (dbExcel) - shifting old B->C and C->D.
$ws.Columns.Item(2).Insert()

# New header text for the inserted column
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query text for the inserted column's data row
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Boxer']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").Value = $statQuery

# The inserted column copies the formatting of the column to its left
# (wrap-text style from A2) - reset it back to the default "Normal" style
# to match the un-styled source cell.
$ws.Range("B2").Style = "Normal"

# Column widths: A/C/D keep their original widths untouched; only the new
# column B needs an explicit (wide, best-fit-like) width.
$ws.Columns.Item(2).ColumnWidth = 254.75

# Reset view: no frozen/scrolled top-left cell, selection on A2
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()
